$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.022.57'
$ws.Range("E2").Value = '  -1.36%  '
$ws.Range("D3").Value = '2.378.93'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.24'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.11'
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").Value = '2.379.17'
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.09'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.79'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("E15").Value = '  -0.98%  '
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").Value = '59.879.98'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Value = '2.376.31'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.13'
$ws.Range("E19").Value = '  +13.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.52'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.53'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -2.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '63.98'
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '556.63'
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.13'
$ws.Range("E28").Value = '  -5.50%  '
$ws.Range("D29").Value = '2.488.25'
$ws.Range("D30").Value = '0.0₃0927'
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.99'
$ws.Range("E31").Value = '  +2.43%  '
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.79'
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.132'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.44'
$ws.Range("E36").Value = '  +5.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.26'
$ws.Range("E37").Value = '  +3.96%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.55'
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.17'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.03'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.60'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  +3.90%  '
$ws.Range("D46").Value = '0.0₆0296'
$ws.Range("E46").Value = '  +5.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.49'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.586'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0500'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.15'
$ws.Range("E51").Value = '  -0.59%  '
